$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 91 (existing rows 91-141 shift down to 92-142).
$ws.Rows("91:91").Insert()

# Populate the newly inserted row 91 with the new price record.
$ws.Range("A91").Value = 10
$ws.Range("B91").Value = "Vega Modelo de Temuco"
$ws.Range("C91").Value = "La Araucanía"
$ws.Range("D91").Value = 44455
$ws.Range("D91").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E91").Value = 9
$ws.Range("F91").Value = 100112039
$ws.Range("G91").Value = "Ciboulette"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 20
$ws.Range("K91").Value = 7000
$ws.Range("L91").Value = 8000
$ws.Range("M91").Value = 7500
$ws.Range("N91").Value = "$/docena de atados"
$ws.Range("O91").Value = "Provincia de Cautín"
$ws.Range("P91").Value = 2500
$ws.Range("Q91").Value = 3
$ws.Range("R91").Value = "Hortaliza"
